$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 268.1
$ws.Range("I5").Value = 268.1
$ws.Range("K5").Value = 268.1
$ws.Range("M5").Value = -153.1
$ws.Range("H40").Value = 7415.875
$ws.Range("I40").Value = 5800
$ws.Range("K40").Value = 5800
$ws.Range("M40").Value = -5625
$ws.Range("H43").Value = 1952.1428
$ws.Range("I43").Value = 1892.5
$ws.Range("K43").Value = 1892.5
$ws.Range("M43").Value = -1823.5
$ws.Range("H98").Value = 3636.074
$ws.Range("I98").Value = 4062.2727
$ws.Range("K98").Value = 4062.2727
$ws.Range("M98").Value = -2564.2727
$ws.Range("H99").Value = 400.26086
$ws.Range("I99").Value = 264.6316
$ws.Range("J99").Value = 1044.5
$ws.Range("K99").Value = 793.8948
$ws.Range("L99").Value = 3133.5
$ws.Range("M99").Value = 704.1052
$ws.Range("N99").Value = -6129.5
$ws.Range("H113").Value = 3139.5833
$ws.Range("I113").Value = 2700
$ws.Range("J113").Value = 3179.5454
$ws.Range("K113").Value = 2700
$ws.Range("L113").Value = 3179.5454
$ws.Range("M113").Value = 554
$ws.Range("N113").Value = -9687.545399999999
$ws.Range("H116").Value = 6380.8423
$ws.Range("I116").Value = 5731.2856
$ws.Range("K116").Value = 5731.2856
$ws.Range("M116").Value = -2289.2856
$ws.Range("H122").Value = 3636.074
$ws.Range("I122").Value = 4062.2727
$ws.Range("K122").Value = 12186.8181
$ws.Range("M122").Value = -9736.8181
$ws.Range("H132").Value = 4855.4326
$ws.Range("I132").Value = 5378.533
$ws.Range("J132").Value = 2613.5715
$ws.Range("K132").Value = 16135.599
$ws.Range("L132").Value = 7840.7145
$ws.Range("M132").Value = -13605.599
$ws.Range("N132").Value = -12900.7145
$ws.Range("H133").Value = 106999
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H135").Value = 1107.7693
$ws.Range("I135").Value = 986.7
$ws.Range("K135").Value = 8880.300000000001
$ws.Range("M135").Value = -6345.300000000001
$ws.Range("H138").Value = 3659.48
$ws.Range("I138").Value = 3099.2
$ws.Range("J138").Value = 4499.9
$ws.Range("K138").Value = 9297.599999999999
$ws.Range("L138").Value = 13499.7
$ws.Range("M138").Value = -4157.599999999999
$ws.Range("N138").Value = -23779.7
$ws.Range("H141").Value = 2092.8235
$ws.Range("I141").Value = 2092.8235
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6278.470499999999
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1935
$ws.Range("I2").Value = 1849.375
$ws.Range("K2").Value = 1849.375
$ws.Range("M2").Value = -1736.375
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H42").Value = 17000
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H45").Value = 24004.63
$ws.Range("I45").Value = 31792
$ws.Range("K45").Value = 31792
$ws.Range("M45").Value = -31415
$ws.Range("H74").Value = 178539.22
$ws.Range("I74").Value = 227065.44
$ws.Range("J74").Value = 5231.2856
$ws.Range("K74").Value = 227065.44
$ws.Range("L74").Value = 5231.2856
$ws.Range("M74").Value = -226191.44
$ws.Range("N74").Value = -6979.2856
$ws.Range("H77").Value = 178539.22
$ws.Range("I77").Value = 227065.44
$ws.Range("J77").Value = 5231.2856
$ws.Range("K77").Value = 1135327.2
$ws.Range("L77").Value = 26156.428
$ws.Range("M77").Value = -1130959.2
$ws.Range("N77").Value = -34892.428
$ws.Range("H95").Value = 39999
$ws.Range("J95").Value = 39999
$ws.Range("L95").Value = 39999
$ws.Range("N95").Value = -45491
$ws.Range("H102").Value = 3196.4736
$ws.Range("I102").Value = 2269
$ws.Range("J102").Value = 6674.5
$ws.Range("K102").Value = 2269
$ws.Range("L102").Value = 6674.5
$ws.Range("M102").Value = -647
$ws.Range("N102").Value = -9918.5
$ws.Range("H116").Value = 1935
$ws.Range("I116").Value = 1849.375
$ws.Range("K116").Value = 1849.375
$ws.Range("M116").Value = 444.625
$ws.Range("H132").Value = 2080.3
$ws.Range("I132").Value = 1774.5862
$ws.Range("K132").Value = 5323.7586
$ws.Range("M132").Value = -2793.7586

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1935
$ws.Range("I3").Value = 1849.375
$ws.Range("K3").Value = 1849.375
$ws.Range("M3").Value = -1735.375
$ws.Range("H22").Value = 597.5263
$ws.Range("I22").Value = 602.94446
$ws.Range("K22").Value = 602.94446
$ws.Range("M22").Value = -429.94446
$ws.Range("H26").Value = 3867.75
$ws.Range("I26").Value = 3867.75
$ws.Range("K26").Value = 3867.75
$ws.Range("M26").Value = -3575.75
$ws.Range("H55").Value = 30777
$ws.Range("J55").Value = 30777
$ws.Range("L55").Value = 30777
$ws.Range("N55").Value = -31323
$ws.Range("H94").Value = 100001610
$ws.Range("I94").Value = 125001220
$ws.Range("K94").Value = 125001220
$ws.Range("M94").Value = -125000769
$ws.Range("H116").Value = 79990
$ws.Range("J116").Value = 79990
$ws.Range("L116").Value = 79990
$ws.Range("N116").Value = -89168
$ws.Range("H134").Value = 4731.22
$ws.Range("I134").Value = 2072.7334
$ws.Range("J134").Value = 13276.357
$ws.Range("K134").Value = 6218.2002
$ws.Range("L134").Value = 39829.071
$ws.Range("M134").Value = -3683.2002
$ws.Range("N134").Value = -44899.071

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1879.3529
$ws.Range("I16").Value = 1901.5555
$ws.Range("K16").Value = 1901.5555
$ws.Range("M16").Value = -1614.5555
$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H41").Value = 55975
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 55975
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -56831
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H58").Value = 2214.5925
$ws.Range("I58").Value = 1561.5
$ws.Range("K58").Value = 1561.5
$ws.Range("M58").Value = -1358.5
$ws.Range("H63").Value = 70000
$ws.Range("J63").Value = 70000
$ws.Range("L63").Value = 70000
$ws.Range("N63").Value = -71372
$ws.Range("H66").Value = 70000
$ws.Range("J66").Value = 70000
$ws.Range("L66").Value = 210000
$ws.Range("N66").Value = -216864
$ws.Range("H86").Value = 9985
$ws.Range("J86").Value = 14998.75
$ws.Range("L86").Value = 14998.75
$ws.Range("N86").Value = -17244.75
$ws.Range("H89").Value = 9985
$ws.Range("J89").Value = 14998.75
$ws.Range("L89").Value = 74993.75
$ws.Range("N89").Value = -86225.75
$ws.Range("H94").Value = 1637.8572
$ws.Range("J94").Value = 1978.25
$ws.Range("L94").Value = 1978.25
$ws.Range("N94").Value = -2880.25
$ws.Range("H99").Value = 6170.875
$ws.Range("J99").Value = 7492.6665
$ws.Range("L99").Value = 7492.6665
$ws.Range("N99").Value = -10488.6665
$ws.Range("H102").Value = 29309.5
$ws.Range("I102").Value = 30219
$ws.Range("J102").Value = 28400
$ws.Range("K102").Value = 30219
$ws.Range("L102").Value = 28400
$ws.Range("M102").Value = -27785
$ws.Range("N102").Value = -33268
$ws.Range("H105").Value = 2256
$ws.Range("I105").Value = 1669.4546
$ws.Range("J105").Value = 3062.5
$ws.Range("K105").Value = 1669.4546
$ws.Range("L105").Value = 3062.5
$ws.Range("M105").Value = 77.54539999999997
$ws.Range("N105").Value = -6556.5
$ws.Range("H107").Value = 2778431
$ws.Range("I107").Value = 4545794.5
$ws.Range("J107").Value = 1145.1428
$ws.Range("K107").Value = 4545794.5
$ws.Range("L107").Value = 1145.1428
$ws.Range("M107").Value = -4543874.5
$ws.Range("N107").Value = -4985.1428
$ws.Range("H113").Value = 1879.3529
$ws.Range("I113").Value = 1901.5555
$ws.Range("K113").Value = 1901.5555
$ws.Range("M113").Value = 268.4445000000001
$ws.Range("H120").Value = 49998.5
$ws.Range("J120").Value = 49998.5
$ws.Range("L120").Value = 49998.5
$ws.Range("N120").Value = -57256.5
$ws.Range("H121").Value = 49997
$ws.Range("J121").Value = 49997
$ws.Range("L121").Value = 49997
$ws.Range("N121").Value = -52617
$ws.Range("H122").Value = 3153
$ws.Range("I122").Value = 3241.8572
$ws.Range("J122").Value = 2842
$ws.Range("K122").Value = 9725.571599999999
$ws.Range("L122").Value = 8526
$ws.Range("M122").Value = -7275.571599999999
$ws.Range("N122").Value = -13426
$ws.Range("H126").Value = 6170.875
$ws.Range("J126").Value = 7492.6665
$ws.Range("L126").Value = 22477.9995
$ws.Range("N126").Value = -27417.9995
$ws.Range("H133").Value = 71138.39999999999
$ws.Range("J133").Value = 79423
$ws.Range("L133").Value = 79423
$ws.Range("N133").Value = -84483
$ws.Range("H134").Value = 2272.3774
$ws.Range("I134").Value = 1990.091
$ws.Range("K134").Value = 5970.272999999999
$ws.Range("M134").Value = -3435.272999999999
$ws.Range("H135").Value = 96832.336
$ws.Range("J135").Value = 96832.336
$ws.Range("L135").Value = 96832.336
$ws.Range("N135").Value = -106972.336
$ws.Range("H136").Value = 2214.5925
$ws.Range("I136").Value = 1561.5
$ws.Range("K136").Value = 4684.5
$ws.Range("M136").Value = -2134.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 124.5
$ws.Range("I12").Value = 26.5
$ws.Range("J12").Value = 157.16667
$ws.Range("K12").Value = 79.5
$ws.Range("L12").Value = 471.50001
$ws.Range("M12").Value = 93.5
$ws.Range("N12").Value = -817.50001
$ws.Range("H17").Value = 976.625
$ws.Range("I17").Value = 128.75
$ws.Range("J17").Value = 1824.5
$ws.Range("K17").Value = 386.25
$ws.Range("L17").Value = 5473.5
$ws.Range("M17").Value = -217.25
$ws.Range("N17").Value = -5811.5
$ws.Range("H86").Value = 1985.3077
$ws.Range("J86").Value = 2023.3334
$ws.Range("L86").Value = 6070.0002
$ws.Range("N86").Value = -8442.0002
$ws.Range("H89").Value = 1985.3077
$ws.Range("J89").Value = 2023.3334
$ws.Range("L89").Value = 18210.0006
$ws.Range("N89").Value = -30066.0006
$ws.Range("H132").Value = 3399.7273
$ws.Range("J132").Value = 3399.7273
$ws.Range("L132").Value = 30597.5457
$ws.Range("N132").Value = -35657.5457

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 23818.834
$ws.Range("J32").Value = 23937.5
$ws.Range("L32").Value = 23937.5
$ws.Range("N32").Value = -24529.5
$ws.Range("H97").Value = 1250.3096
$ws.Range("I97").Value = 1009.34375
$ws.Range("J97").Value = 2021.4
$ws.Range("K97").Value = 1009.34375
$ws.Range("L97").Value = 2021.4
$ws.Range("M97").Value = -513.34375
$ws.Range("N97").Value = -3013.4
$ws.Range("H113").Value = 6475.75
$ws.Range("I113").Value = 3055.5715
$ws.Range("J113").Value = 11264
$ws.Range("K113").Value = 3055.5715
$ws.Range("L113").Value = 11264
$ws.Range("M113").Value = -885.5715
$ws.Range("N113").Value = -15604
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 2224.76
$ws.Range("I132").Value = 1870.5
$ws.Range("K132").Value = 5611.5
$ws.Range("M132").Value = -3081.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4299.3335
$ws.Range("I40").Value = 4299.3335
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4299.3335
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H55").Value = 680.16
$ws.Range("I55").Value = 573.4706
$ws.Range("J55").Value = 906.875
$ws.Range("K55").Value = 573.4706
$ws.Range("L55").Value = 906.875
$ws.Range("M55").Value = -400.4706
$ws.Range("N55").Value = -1252.875
$ws.Range("H64").Value = 35999.57
$ws.Range("J64").Value = 35999.57
$ws.Range("L64").Value = 35999.57
$ws.Range("N64").Value = -36449.57
$ws.Range("H67").Value = 35999.57
$ws.Range("J67").Value = 35999.57
$ws.Range("L67").Value = 35999.57
$ws.Range("N67").Value = -37559.57
$ws.Range("H100").Value = 6725.1
$ws.Range("I100").Value = 4600.3335
$ws.Range("J100").Value = 7635.7144
$ws.Range("K100").Value = 4600.3335
$ws.Range("L100").Value = 7635.7144
$ws.Range("M100").Value = -4059.3335
$ws.Range("N100").Value = -8717.714400000001
$ws.Range("H122").Value = 6889.1304
$ws.Range("I122").Value = 4620.5
$ws.Range("K122").Value = 13861.5
$ws.Range("M122").Value = -11411.5
$ws.Range("H136").Value = 3204.0908
$ws.Range("I136").Value = 3174.75
$ws.Range("K136").Value = 9524.25
$ws.Range("M136").Value = -6974.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4671.5557
$ws.Range("I81").Value = 4109.6
$ws.Range("K81").Value = 8219.200000000001
$ws.Range("M81").Value = -7158.200000000001
$ws.Range("H84").Value = 4671.5557
$ws.Range("I84").Value = 4109.6
$ws.Range("K84").Value = 41096
$ws.Range("M84").Value = -35792
$ws.Range("H96").Value = 5198
$ws.Range("I96").Value = 2999.75
$ws.Range("K96").Value = 2999.75
$ws.Range("M96").Value = -1626.75
$ws.Range("H122").Value = 41668496
$ws.Range("I122").Value = 2196.8
$ws.Range("K122").Value = 6590.400000000001
$ws.Range("M122").Value = -4140.400000000001
$ws.Range("H132").Value = 3138.3416
$ws.Range("I132").Value = 2807.3157
$ws.Range("K132").Value = 8421.947100000001
$ws.Range("M132").Value = -5891.947100000001
$ws.Range("H136").Value = 248971.33
$ws.Range("I136").Value = 4214.564
$ws.Range("J136").Value = 983241.6
$ws.Range("K136").Value = 12643.692
$ws.Range("L136").Value = 2949724.8
$ws.Range("M136").Value = -10093.692
$ws.Range("N136").Value = -2954824.8
